$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at position 483, pushing existing rows 483-552 down to 486-555.
$ws.Range("A483:A485").EntireRow.Insert()

# Populate new row 483
$ws.Range("A483").Value2 = 7
$ws.Range("B483").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C483").Value2 = "Ñuble"
$ws.Range("D483").Value2 = 45127
$ws.Range("E483").Value2 = 16
$ws.Range("F483").Value2 = 100112023
$ws.Range("G483").Value2 = "Brócoli"
$ws.Range("H483").Value2 = "Sin especificar"
$ws.Range("I483").Value2 = "Primera"
$ws.Range("J483").Value2 = 300
$ws.Range("K483").Value2 = 900
$ws.Range("L483").Value2 = 900
$ws.Range("M483").Value2 = 900
$ws.Range("N483").Value2 = "$/unidad"
$ws.Range("O483").Value2 = "Provincia de Diguillín"
$ws.Range("P483").Value2 = 900
$ws.Range("Q483").Value2 = 1
$ws.Range("R483").Value2 = "Hortaliza"

# Populate new row 484
$ws.Range("A484").Value2 = 7
$ws.Range("B484").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C484").Value2 = "Ñuble"
$ws.Range("D484").Value2 = 45127
$ws.Range("E484").Value2 = 16
$ws.Range("F484").Value2 = 100112023
$ws.Range("G484").Value2 = "Brócoli"
$ws.Range("H484").Value2 = "Sin especificar"
$ws.Range("I484").Value2 = "Primera"
$ws.Range("J484").Value2 = 250
$ws.Range("K484").Value2 = 1000
$ws.Range("L484").Value2 = 1000
$ws.Range("M484").Value2 = 1000
$ws.Range("N484").Value2 = "$/unidad"
$ws.Range("O484").Value2 = "Región del Maule"
$ws.Range("P484").Value2 = 1000
$ws.Range("Q484").Value2 = 1
$ws.Range("R484").Value2 = "Hortaliza"

# Populate new row 485
$ws.Range("A485").Value2 = 7
$ws.Range("B485").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C485").Value2 = "Ñuble"
$ws.Range("D485").Value2 = 45127
$ws.Range("E485").Value2 = 16
$ws.Range("F485").Value2 = 100112023
$ws.Range("G485").Value2 = "Brócoli"
$ws.Range("H485").Value2 = "Sin especificar"
$ws.Range("I485").Value2 = "Segunda"
$ws.Range("J485").Value2 = 300
$ws.Range("K485").Value2 = 800
$ws.Range("L485").Value2 = 800
$ws.Range("M485").Value2 = 800
$ws.Range("N485").Value2 = "$/unidad"
$ws.Range("O485").Value2 = "Región del Maule"
$ws.Range("P485").Value2 = 800
$ws.Range("Q485").Value2 = 1
$ws.Range("R485").Value2 = "Hortaliza"
